$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "LOM3206 - Eletrônica (Requisito)" requirement row entirely.
# This shifts rows 24-26 up to 23-25 and keeps all other formatting intact.
$ws.Range("A23").EntireRow.Delete()

# Update the remaining requisito rows' text to reflect the new requirement list.
$ws.Range("B24").Value = "LOM3234 -  Óptica Física  (Requisito)`n"
$ws.Range("C24").Value = "LOM3234 -  Óptica Física  (Requisito)`n"

$ws.Range("B25").Value = "LOM3263 -  Eletrônica Fundamental e Aplicada  (Requisito)`n"
$ws.Range("C25").Value = "LOM3263 -  Eletrônica Fundamental e Aplicada  (Requisito)`n"
